$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(7).Copy() | Out-Null
$ws.Rows.Item(1).Insert() | Out-Null

$ws.Range("C12").Select() | Out-Null
